$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 'mug'
    3 = 'time'
    4 = 'capital'
    5 = 'rage'
    6 = 'navy'
    7 = 'climb'
    8 = 'boat'
    9 = 'measuring'
    10 = 'ire'
    11 = 'horn'
    12 = 'fury'
    13 = 'sandpaper'
    14 = 'honey'
    15 = 'green'
    16 = 'country'
    17 = 'soft'
    18 = 'united states'
    19 = 'uneven'
    20 = 'swim'
    21 = 'slaw'
    22 = 'cautious'
    23 = 'black'
    24 = 'palace'
    25 = 'lungs'
    26 = 'chocolate'
    27 = 'legs'
    28 = 'leg'
    29 = 'danger'
    30 = 'happiness'
    31 = 'road'
    32 = 'room'
    33 = 'desire'
    34 = 'bad'
    35 = 'molehill'
    36 = 'clouds'
    37 = 'disease'
    38 = 'melody'
    39 = 'water'
    40 = 'hand'
    41 = 'glass'
    42 = 'throne'
    43 = 'house'
    44 = 'shout'
    45 = 'pollution'
    46 = 'tobacco'
    47 = 'gloves'
    48 = 'eat'
    49 = 'picture'
    50 = 'dump'
    51 = 'dentist'
    52 = 'wood'
    53 = 'flow'
    54 = 'seat'
    55 = 'elastic'
    56 = 'snore'
    57 = 'keys'
    58 = 'slice'
    59 = 'see'
    60 = 'stein'
    61 = 'vines'
    62 = 'veal'
    63 = 'head'
    64 = 'sing'
    65 = 'lamb'
    66 = 'highway'
    67 = 'quick'
    68 = 'tune'
    69 = 'temper'
    70 = 'saucer'
    71 = 'emblem'
    72 = 'blow'
    73 = 'lake'
    74 = 'halt'
    75 = 'anthem'
    76 = 'fair'
    77 = 'goat'
    78 = 'automobile'
    79 = 'lawyer'
    80 = 'scraps'
    81 = 'summer'
    82 = 'heat'
    83 = 'peak'
    84 = 'quill'
    85 = 'leak'
    86 = 'cook'
    87 = 'fragrance'
    88 = 'marker'
    89 = 'crowded'
    90 = 'wake'
    91 = 'bike'
    92 = 'swivel'
    93 = 'cushion'
    94 = 'beer'
    95 = 'thin'
    96 = 'mad'
    97 = 'sour'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 1).Value = $values[$row]
}
